$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates: force text format to preserve exact string
# representation (these values look numeric and would otherwise be
# auto-converted to floating point numbers, losing formatting such as
# trailing zeros or multi-dot separators).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.898.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.625.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0611"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0914"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.858.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.633.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.958.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0707"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0490"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.425.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.559"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0502"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "54.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.996"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.766.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "88.82"
$ws.Range("D50").Style = "Normal"

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  +10.74%  "
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("E14").Value = "  +6.47%  "
$ws.Range("E15").Value = "  +5.09%  "
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("E17").Value = "  +16.68%  "
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  +3.53%  "
$ws.Range("E23").Value = "  +4.55%  "
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("E30").Value = "  +3.14%  "
$ws.Range("E31").Value = "  +5.47%  "
$ws.Range("E32").Value = "  +4.17%  "
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").Value = "  +7.02%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  +4.27%  "
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("E45").Value = "  +5.35%  "
$ws.Range("E46").Value = "  +17.02%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("E51").Value = "  +3.90%  "
